# Auto-generated cell updates for cryptos worksheet refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.820.36"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.542.94"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "'205.93"
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("D9").Value = "'21.41"
$ws.Range("E9").Value = "  -2.62%  "
$ws.Range("D11").Value = "'0.0854"
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("D12").Value = "1.762.98"
$ws.Range("E12").Value = "  -1.49%  "
$ws.Range("D13").Value = "1.546.48"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("E14").Value = "  -1.53%  "
$ws.Range("D15").Value = "'0.509"
$ws.Range("E15").Value = "  -1.03%  "
$ws.Range("D16").Value = "26.824.48"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "'61.25"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").Value = "'214.74"
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("D19").Value = "'7.23"
$ws.Range("E19").Value = "  -2.49%  "
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("D22").Value = "'3.99"
$ws.Range("E22").Value = "  -3.06%  "
$ws.Range("E23").Value = "  -1.25%  "
$ws.Range("E24").Value = "  -2.91%  "
$ws.Range("D25").Value = "'152.73"
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("D26").Value = "'6.60"
$ws.Range("E26").Value = "  -2.06%  "
$ws.Range("D27").Value = "'14.84"
$ws.Range("E27").Value = "  -0.87%  "
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("D29").Value = "'0.103"
$ws.Range("E29").Value = "  -0.72%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.0458"
$ws.Range("E30").Value = "  -1.92%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.10"
$ws.Range("E31").Value = "  -1.63%  "
$ws.Range("E32").Value = "  +1.65%  "
$ws.Range("D33").Value = "1.367.37"
$ws.Range("E33").Value = "  -2.08%  "
$ws.Range("E34").Value = "  +0.48%  "
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("D36").Value = "'0.964"
$ws.Range("E36").Value = "  +2.94%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("E38").Value = "  +1.13%  "
$ws.Range("D39").Value = "'0.520"
$ws.Range("E39").Value = "  -1.55%  "
$ws.Range("D40").Value = "'5.79"
$ws.Range("E40").Value = "  +9.04%  "
$ws.Range("E41").Value = "  -1.08%  "
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("D43").Value = "'0.990"
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("E44").Value = "  +1.37%  "
$ws.Range("D45").Value = "'63.17"
$ws.Range("E45").Value = "  -0.20%  "
$ws.Range("E46").Value = "  -3.24%  "
$ws.Range("D47").Value = "1.677.34"
$ws.Range("E47").Value = "  -1.50%  "
$ws.Range("D48").Value = "'84.15"
$ws.Range("E48").Value = "  -2.28%  "
$ws.Range("D49").Value = "'0.0510"
$ws.Range("E49").Value = "  +3.62%  "
$ws.Range("D50").Value = "0.0₇0969"
$ws.Range("E50").Value = "  -1.70%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  +0.19%  "
